$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price (D) column cells are treated as plain text so values
# like "1.00", "7.20", "2.90", "397.10" keep their exact formatting instead
# of being auto-converted to numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.445.15"
$ws.Range("E2").Value = "  -1.01%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.761.31"
$ws.Range("E3").Value = "  -2.16%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.73"
$ws.Range("E5").Value = "  -1.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.93"
$ws.Range("E6").Value = "  +0.95%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.761.18"
$ws.Range("E7").Value = "  -2.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("E9").Value = "  -0.37%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.163"
$ws.Range("E10").Value = "  -0.36%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.47"
$ws.Range("E11").Value = "  +0.57%  "

$ws.Range("E12").Value = "  -0.94%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000271"
$ws.Range("E13").Value = "  +3.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.56"
$ws.Range("E14").Value = "  -1.00%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.395.33"
$ws.Range("E15").Value = "  -2.17%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.759.21"
$ws.Range("E16").Value = "  -2.29%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.84"
$ws.Range("E17").Value = "  +3.87%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.491.75"
$ws.Range("E18").Value = "  -1.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.20"
$ws.Range("E19").Value = "  -2.45%  "

$ws.Range("E20").Value = "  +0.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.57"
$ws.Range("E21").Value = "  -3.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "467.59"
$ws.Range("E22").Value = "  +0.27%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.722"
$ws.Range("E23").Value = "  -1.30%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.23"
$ws.Range("E26").Value = "  -0.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.15"
$ws.Range("E27").Value = "  +0.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.27"
$ws.Range("E28").Value = "  +2.52%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.90"
$ws.Range("E30").Value = "  -2.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.908.07"
$ws.Range("E31").Value = "  -2.21%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.64"
$ws.Range("E32").Value = "  +0.58%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.25"
$ws.Range("E33").Value = "  -2.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.35"
$ws.Range("E34").Value = "  -2.42%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.14"
$ws.Range("E35").Value = "  -4.79%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.728.48"
$ws.Range("E36").Value = "  -2.12%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.84"
$ws.Range("E37").Value = "  +6.95%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.105"
$ws.Range("E38").Value = "  +0.37%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.996"
$ws.Range("E41").Value = "  -2.62%  "

$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.313"
$ws.Range("E43").Value = "  +0.11%  "

$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.69"
$ws.Range("E45").Value = "  +0.83%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.95"
$ws.Range("E46").Value = "  -1.21%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.76"
$ws.Range("E47").Value = "  -2.51%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "397.10"
$ws.Range("E48").Value = "  -5.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.68"
$ws.Range("E49").Value = "  -0.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000267"
$ws.Range("E50").Value = "  -10.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0354"
$ws.Range("E51").Value = "  -0.91%  "

$ws.Range("B24").Value = "PEPE"
$ws.Range("C24").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000147"
$ws.Range("E24").Value = "  -7.79%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.71"
$ws.Range("E25").Value = "  +0.96%  "

$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.88"
$ws.Range("E39").Value = "  -0.32%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.137"
$ws.Range("E40").Value = "  -2.39%  "

